$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $style = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = $style
}

Set-TextValue $ws.Range("D2") "26.380.23"
Set-TextValue $ws.Range("E2") "  -3.95%  "

Set-TextValue $ws.Range("D3") "1.766.16"
Set-TextValue $ws.Range("E3") "  -3.23%  "

Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.10%  "

Set-TextValue $ws.Range("D5") "1.002"
Set-TextValue $ws.Range("E5") "  +0.03%  "

Set-TextValue $ws.Range("D6") "305.62"
Set-TextValue $ws.Range("E6") "  -2.27%  "

Set-TextValue $ws.Range("D7") "0.4298"
Set-TextValue $ws.Range("E7") "  +1.23%  "

Set-TextValue $ws.Range("D8") "0.3632"
Set-TextValue $ws.Range("E8") "  +0.38%  "

Set-TextValue $ws.Range("D9") "0.07077"
Set-TextValue $ws.Range("E9") "  -1.89%  "

Set-TextValue $ws.Range("D10") "0.8446"
Set-TextValue $ws.Range("E10") "  -1.84%  "

Set-TextValue $ws.Range("D11") "20.25"
Set-TextValue $ws.Range("E11") "  -1.75%  "

Set-TextValue $ws.Range("D12") "1.753.77"
Set-TextValue $ws.Range("E12") "  -4.05%  "

Set-TextValue $ws.Range("D13") "5.245"
Set-TextValue $ws.Range("E13") "  -2.75%  "

Set-TextValue $ws.Range("D14") "6.427"
Set-TextValue $ws.Range("E14") "  -0.81%  "

Set-TextValue $ws.Range("D15") "0.06801"
Set-TextValue $ws.Range("E15") "  -1.83%  "

Set-TextValue $ws.Range("E16") "  +0.19%  "

Set-TextValue $ws.Range("D17") "79.18"
Set-TextValue $ws.Range("E17") "  -1.46%  "

Set-TextValue $ws.Range("D18") "0.000008599"
Set-TextValue $ws.Range("E18") "  -3.18%  "

Set-TextValue $ws.Range("D19") "1.001"
Set-TextValue $ws.Range("E19") "  -0.06%  "

Set-TextValue $ws.Range("D20") "15.04"
Set-TextValue $ws.Range("E20") "  -1.95%  "

Set-TextValue $ws.Range("D21") "26.389.08"
Set-TextValue $ws.Range("E21") "  -4.34%  "

Set-TextValue $ws.Range("E22") "  -1.79%  "

Set-TextValue $ws.Range("E23") "  +2.50%  "

Set-TextValue $ws.Range("D24") "2.014.95"
Set-TextValue $ws.Range("E24") "  -2.09%  "

Set-TextValue $ws.Range("D25") "152.72"
Set-TextValue $ws.Range("E25") "  -1.48%  "

Set-TextValue $ws.Range("D26") "1.846"
Set-TextValue $ws.Range("E26") "  -7.06%  "

Set-TextValue $ws.Range("D27") "18.15"
Set-TextValue $ws.Range("E27") "  -2.90%  "

Set-TextValue $ws.Range("D28") "5.077"
Set-TextValue $ws.Range("E28") "  -1.38%  "

Set-TextValue $ws.Range("D29") "114.07"
Set-TextValue $ws.Range("E29") "  -0.13%  "

Set-TextValue $ws.Range("D30") "1.703"
Set-TextValue $ws.Range("E30") "  -5.22%  "

Set-TextValue $ws.Range("D31") "0.08929"
Set-TextValue $ws.Range("E31") "  +0.87%  "

Set-TextValue $ws.Range("D32") "0.7304"
Set-TextValue $ws.Range("E32") "  -2.21%  "

Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "4.334"
Set-TextValue $ws.Range("E33") "  -4.41%  "

Set-TextValue $ws.Range("B34") "ARBITRUM"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "1.112"
Set-TextValue $ws.Range("E34") "  -0.63%  "

Set-TextValue $ws.Range("D35") "2.763"
Set-TextValue $ws.Range("E35") "  -7.22%  "

Set-TextValue $ws.Range("E36") "  +0.02%  "

Set-TextValue $ws.Range("D37") "1.075"
Set-TextValue $ws.Range("E37") "  -0.93%  "

Set-TextValue $ws.Range("D38") "0.05127"
Set-TextValue $ws.Range("E38") "  -2.99%  "

Set-TextValue $ws.Range("D39") "0.01892"
Set-TextValue $ws.Range("E39") "  -1.32%  "

Set-TextValue $ws.Range("D40") "0.1610"
Set-TextValue $ws.Range("E40") "  -1.97%  "

Set-TextValue $ws.Range("D41") "0.4908"
Set-TextValue $ws.Range("E41") "  -3.08%  "

Set-TextValue $ws.Range("D42") "2.518"
Set-TextValue $ws.Range("E42") "  -9.36%  "

Set-TextValue $ws.Range("D43") "6.211"
Set-TextValue $ws.Range("E43") "  -3.49%  "

Set-TextValue $ws.Range("E44") "  -3.54%  "

Set-TextValue $ws.Range("D45") "104.93"
Set-TextValue $ws.Range("E45") "  -0.66%  "

Set-TextValue $ws.Range("E46") "  +0.02%  "

Set-TextValue $ws.Range("D47") "10.11"
Set-TextValue $ws.Range("E47") "  -2.70%  "

Set-TextValue $ws.Range("D48") "0.06187"
Set-TextValue $ws.Range("E48") "  -4.01%  "

Set-TextValue $ws.Range("D49") "0.4483"
Set-TextValue $ws.Range("E49") "  -4.13%  "

Set-TextValue $ws.Range("D50") "1.577"
Set-TextValue $ws.Range("E50") "  -2.11%  "

Set-TextValue $ws.Range("D51") "1.744"
Set-TextValue $ws.Range("E51") "  +2.93%  "
